$wb = $excel.ActiveWorkbook

# Updates to the "想去人数" (want-to-go count) column F for sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    4  = 66
    5  = 240
    7  = 70
    10 = 48
    11 = 32
    12 = 101
    13 = 2263
    15 = 33
    16 = 516
    17 = 524
    18 = 159
    19 = 81
    21 = 47
    22 = 1745
    23 = 3886
    27 = 1159
    28 = 218
    29 = 2053
    30 = 568
    33 = 283
    34 = 416
    36 = 677
    38 = 402
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Updates to the "想去人数" (want-to-go count) column F for sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    4  = 66
    5  = 240
    7  = 70
    10 = 48
    11 = 32
    12 = 101
    13 = 2263
    16 = 33
    17 = 516
    18 = 524
    19 = 159
    20 = 81
    22 = 47
    23 = 1745
    24 = 3887
    28 = 1159
    29 = 218
    30 = 2053
    31 = 568
    34 = 283
    35 = 416
    37 = 677
    39 = 402
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
